$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 2.4
$ws.Range("G2").Value = 2.48
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 3.8
$ws.Range("J2").Value = 3.15
$ws.Range("K2").Value = 3.3
$ws.Range("L2").Value = 1.49
$ws.Range("M2").Value = 1.09
$ws.Range("P2").Value = 1.69
$ws.Range("W2").Value = 1.56
$ws.Range("Z2").Value = 24

# Row 3
$ws.Range("F3").Value = 1.75
$ws.Range("G3").Value = 1.89
$ws.Range("H3").Value = 5.3
$ws.Range("L3").Value = 1.57
$ws.Range("P3").Value = 1.61
$ws.Range("Q3").Value = 2.14
$ws.Range("S3").Value = 4.8
$ws.Range("Y3").Value = 980
$ws.Range("AJ3").Value = 25
$ws.Range("AN3").Value = 22

# Row 4
$ws.Range("F4").Value = 2.84
$ws.Range("L4").Value = 1.41

# Row 5
$ws.Range("J5").Value = 3.7
$ws.Range("AB5").Value = 7.4
$ws.Range("AC5").Value = 9.6
$ws.Range("AF5").Value = 9.4
$ws.Range("AG5").Value = 11
$ws.Range("AJ5").Value = 18.5
$ws.Range("AL5").Value = 55

# Row 6
$ws.Range("H6").Value = 5.1
$ws.Range("L6").Value = 1.45
